# The "Förändrad" (Changed) date column (C) was bumped by one day
# (2023-09-10 -> 2023-09-11, Excel serial 45179 -> 45180) for every
# data row (rows 2 through 79) on the single worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C79").Value = 45180
